$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Drop the L1_agg_fuel / L2_CEDS_fuel / L3_agg_sector (always-empty) columns and the
# L4_CEDS_sector column (held the "x" flag). Deleting G:J shifts the trailing
# override_normalization / start_continuity / end_continuity columns left into G:I.
$ws2.Range("G1:J1").EntireColumn.Delete() | Out-Null

# New trailing column header that replaces the removed level-specification columns.
$ws2.Range("J1").Value = "user_pct_breakdowns"

# Approximate the bestFit width Excel computed for the new I column ("end_continuity").
$ws2.Columns.Item(9).ColumnWidth = 12.5

# Restore the active-cell selection shown in the sheet (previously L8, now J8 after
# the column shift).
$ws2.Range("J8").Select() | Out-Null
